# Apply the "New crime data collected" edit to the 9th Precinct CompStat weekly report.
#
# - Updates the report header (issue/volume number + the week-covering date range).
# - Updates the Crime Complaints table (rows 14-29) with the refreshed weekly/28-day/
#   year-to-date counts and their derived percent-change columns. A few cells flip
#   between a numeric value and the sheet's "no data" text placeholders ("0" / "***.*");
#   for those we copy the exact cell style from an untouched reference cell (row 30,
#   which this edit never touches) before writing the new value, so the cell keeps the
#   same look (alignment/number format) the source report uses for that placeholder.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header: issue number + reporting week dates (rich-text cells, edit just the run) ---
$ws.Range("A8").Characters(21, 2).Text = "25"
$ws.Range("C9").Characters(27, 9).Text = "6/19/2023"
$ws.Range("C9").Characters(47, 9).Text = "6/25/2023"

# --- Reference cells (stable, untouched by this edit) used to copy exact cell styles
#     when a cell switches between numeric and text ("0" / "***.*") representation. ---
$refText0 = $ws.Range("C30")    # style: right-aligned text, used for the "0" placeholder
$refTextStar = $ws.Range("E30") # style: right-aligned text, used for the "***.*" placeholder
$refNum15 = $ws.Range("J30")    # style: right-aligned integer number
$refNum16 = $ws.Range("K30")    # style: right-aligned 1-decimal percent-change number

# --- Row 14 ---
$refText0.Copy($ws.Range("F14"))

# --- Row 15 ---
$refText0.Copy($ws.Range("C15"))
$ws.Range("M15").Value = 80
$ws.Range("N15").Value = -62.5

# --- Row 16 ---
$ws.Range("C16").Value = 4
$ws.Range("E16").Value = -33.333333333333
$ws.Range("F16").Value = 12
$ws.Range("H16").Value = -25
$ws.Range("I16").Value = 85
$ws.Range("J16").Value = 123
$ws.Range("K16").Value = -30.894308943089
$ws.Range("L16").Value = 26.865671641791
$ws.Range("M16").Value = 10.38961038961
$ws.Range("N16").Value = -78.426395939086

# --- Row 17 ---
$ws.Range("C17").Value = 10
$ws.Range("E17").Value = 400
$ws.Range("F17").Value = 27
$ws.Range("G17").Value = 14
$ws.Range("H17").Value = 92.857142857142
$ws.Range("I17").Value = 126
$ws.Range("J17").Value = 90
$ws.Range("K17").Value = 40
$ws.Range("L17").Value = 65.78947368421
$ws.Range("M17").Value = 61.538461538461
$ws.Range("N17").Value = -53.676470588235

# --- Row 18 ---
$ws.Range("D18").Value = 5
$ws.Range("E18").Value = 20
$ws.Range("F18").Value = 26
$ws.Range("G18").Value = 38
$ws.Range("H18").Value = -31.578947368421
$ws.Range("I18").Value = 138
$ws.Range("J18").Value = 186
$ws.Range("K18").Value = -25.806451612903
$ws.Range("L18").Value = -8
$ws.Range("M18").Value = 25.454545454545
$ws.Range("N18").Value = -63.96866840731

# --- Row 19 ---
$ws.Range("C19").Value = 18
$ws.Range("E19").Value = -14.285714285714
$ws.Range("F19").Value = 85
$ws.Range("G19").Value = 86
$ws.Range("H19").Value = -1.162790697674
$ws.Range("I19").Value = 488
$ws.Range("J19").Value = 464
$ws.Range("K19").Value = 5.172413793103
$ws.Range("L19").Value = 73.665480427046
$ws.Range("M19").Value = 34.806629834254
$ws.Range("N19").Value = -30.681818181818

# --- Row 20 ---
$ws.Range("C20").Value = 2
$ws.Range("F20").Value = 6
$ws.Range("G20").Value = 6
$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 21
$ws.Range("K20").Value = -4.545454545454
$ws.Range("L20").Value = -47.5
$ws.Range("M20").Value = -4.545454545454
$ws.Range("N20").Value = -91.463414634146

# --- Row 21 ---
$ws.Range("C21").Value = 40
$ws.Range("D21").Value = 34
$ws.Range("E21").Value = 17.647058823529
$ws.Range("F21").Value = 157
$ws.Range("G21").Value = 161
$ws.Range("H21").Value = -2.484472049689
$ws.Range("I21").Value = 868
$ws.Range("J21").Value = 900
$ws.Range("K21").Value = -3.555555555555
$ws.Range("L21").Value = 38.658146964856
$ws.Range("M21").Value = 32.317073170731
$ws.Range("N21").Value = -57.199211045364

# --- Row 22 ---
$refText0.Copy($ws.Range("D22"))
$refTextStar.Copy($ws.Range("E22"))

# --- Row 23 ---
$ws.Range("C23").Value = 2
$ws.Range("E23").Value = 0
$ws.Range("F23").Value = 10
$ws.Range("G23").Value = 9
$ws.Range("H23").Value = 11.111111111111
$ws.Range("I23").Value = 67
$ws.Range("J23").Value = 69
$ws.Range("K23").Value = -2.898550724637
$ws.Range("L23").Value = -26.373626373626
$ws.Range("M23").Value = 13.559322033898

# --- Row 24 ---
$ws.Range("C24").Value = 32
$ws.Range("D24").Value = 51
$ws.Range("E24").Value = -37.254901960784
$ws.Range("F24").Value = 111
$ws.Range("G24").Value = 192
$ws.Range("H24").Value = -42.1875
$ws.Range("I24").Value = 699
$ws.Range("J24").Value = 1086
$ws.Range("K24").Value = -35.635359116022
$ws.Range("L24").Value = 61.805555555555
$ws.Range("M24").Value = -10.841836734693

# --- Row 25 ---
$ws.Range("D25").Value = 5
$ws.Range("E25").Value = 80
$ws.Range("I25").Value = 228
$ws.Range("J25").Value = 224
$ws.Range("K25").Value = 1.785714285714
$ws.Range("L25").Value = 49.019607843137
$ws.Range("M25").Value = 3.636363636363

# --- Row 26 ---
$refText0.Copy($ws.Range("C26"))
$ws.Range("F26").Value = 1
$ws.Range("H26").Value = 0

# --- Row 27 ---
$refNum15.Copy($ws.Range("C27"))
$ws.Range("C27").Value = 4
$refNum15.Copy($ws.Range("D27"))
$ws.Range("D27").Value = 3
$refNum16.Copy($ws.Range("E27"))
$ws.Range("E27").Value = 33.333333333333
$ws.Range("F27").Value = 5
$ws.Range("H27").Value = -44.444444444444
$ws.Range("I27").Value = 26
$ws.Range("J27").Value = 46
$ws.Range("K27").Value = -43.478260869565
$ws.Range("L27").Value = -3.703703703703

# --- Row 28 ---
$refText0.Copy($ws.Range("F28"))
$ws.Range("N28").Value = -81.25

# --- Row 29 ---
$refText0.Copy($ws.Range("F29"))
$ws.Range("N29").Value = -75
